$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column headers (shared string text changes)
$ws.Range("D1").Value = "Vahini TCL leakage  power (W)"
$ws.Range("C1").Value = "Vahini TCL power (W)"

# Update some data values to match refreshed calc results
$ws.Range("D2").Value = 0.0000000006336271
$ws.Range("C3").Value = 0.0000378271126
$ws.Range("D3").Value = 0.000000004009252
$ws.Range("C7").Value = 0.00000000000395660816
$ws.Range("D8").Value = 0.00000000001303124

# Update column D width and selected cell to match saved view state
$ws.Columns("D").ColumnWidth = 27
$ws.Range("A11").Select()
